$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: subject "OOPS" already present in A2; only the teacher name changes.
$ws.Range("B2").Value = "Prakhar Golchha"

# Rows 3-6: straightforward row-by-row entry (subject then teacher).
$ws.Range("A3").Value = "CSA "
$ws.Range("B3").Value = "A.K shrivastav"

$ws.Range("A4").Value = "DAA"
$ws.Range("B4").Value = "Akash Sonkar"

$ws.Range("A5").Value = "ES"
$ws.Range("B5").Value = "Devika Sahu"

$ws.Range("A6").Value = "DM"
$ws.Range("B6").Value = "Diwan"

# Rows 7-14: the lab/elective block was entered as a whole subject column
# first, then the teacher column was filled in afterwards.
$ws.Range("A7").Value = "CSA LAB"
$ws.Range("A8").Value = "CSA LAB"
$ws.Range("A9").Value = "OOPS LAB"
$ws.Range("A10").Value = "OPPS LAB"
$ws.Range("A11").Value = "DBMS LAB"
$ws.Range("A12").Value = "DBMS LAB"
$ws.Range("A13").Value = "VL"
$ws.Range("A14").Value = "VL"

$ws.Range("B7").Value = "Anjum khan "
$ws.Range("B8").Value = "Akash Sonkar"
$ws.Range("B9").Value = "Prakhar Golchha"
$ws.Range("B10").Value = "Suyash Sahu"
$ws.Range("B11").Value = "Prakhar Golchha"
$ws.Range("B12").Value = "Palak Lunia"
$ws.Range("B13").Value = "Devika Sahu"
$ws.Range("B14").Value = "Ruchi Agrawal"

# Widen column B to fit the longer teacher names, matching the final layout.
$ws.Columns.Item(2).ColumnWidth = 16.6

# Leave the selection where the last entry was made.
[void]$ws.Range("B14").Select()
